$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy header style (bold, centered, bordered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Fill season record values (Wins, Losses, Ties) for every data row (2-40)
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 29).Value = 84
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 1
}
